$d = $word.ActiveDocument

# Pairs of (old text, new text) for each table-cell answer that was regenerated.
$replacements = @(
    @("227÷6=37, 5",   "718÷2=359, 0"),
    @("923÷6=153, 5",  "350÷7=50, 0"),
    @("731÷4=182, 3",  "394÷8=49, 2"),
    @("145÷7=20, 5",   "219÷4=54, 3"),
    @("369÷2=184, 1",  "386÷6=64, 2"),
    @("373÷3=124, 1",  "998÷8=124, 6"),
    @("475÷5=95, 0",   "616÷8=77, 0"),
    @("622÷3=207, 1",  "829÷3=276, 1"),
    @("395÷3=131, 2",  "659÷5=131, 4"),
    @("509÷2=254, 1",  "661÷7=94, 3"),
    @("864÷7=123, 3",  "727÷7=103, 6"),
    @("200÷8=25, 0",   "569÷4=142, 1"),
    @("658÷9=73, 1",   "110÷5=22, 0"),
    @("959÷8=119, 7",  "539÷2=269, 1"),
    @("127÷4=31, 3",   "102÷5=20, 2"),
    @("896÷3=298, 2",  "740÷7=105, 5"),
    @("758÷9=84, 2",   "943÷6=157, 1"),
    @("999÷4=249, 3",  "699÷3=233, 0"),
    @("167÷3=55, 2",   "573÷3=191, 0"),
    @("295÷9=32, 7",   "523÷6=87, 1"),
    @("372÷7=53, 1",   "117÷2=58, 1"),
    @("950÷2=475, 0",  "392÷7=56, 0"),
    @("927÷9=103, 0",  "941÷4=235, 1"),
    @("732÷9=81, 3",   "941÷3=313, 2"),
    @("123÷9=13, 6",   "571÷4=142, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: could not find '$old' to replace with '$new'"
    }
}
